# Adds a red "NULL" label at the end of the linked-list diagrams that are
# built up across slides 6, 7 and 8 (each slide shows one more node being
# appended; the tail pointer's arrow needs a "NULL" label after it).

$p = $ppt.ActivePresentation

$boxes = @(
    @{ Slide = 6; Left = 3596639; Top = 3607023; Width = 3056709; Height = 707886 },
    @{ Slide = 7; Left = 5682153; Top = 3607023; Width = 3056709; Height = 707886 },
    @{ Slide = 8; Left = 7977051; Top = 3659870; Width = 3056709; Height = 707886 }
)

foreach ($b in $boxes) {
    $s = $p.Slides.Item($b.Slide)
    $tb = $s.Shapes.AddTextbox(1, $b.Left, $b.Top, $b.Width, $b.Height)

    $tb.TextFrame.WordWrap = $true
    $tb.TextFrame.TextRange.Text = "NULL"

    $font = $tb.TextFrame.TextRange.Font
    $font.Size = 40
    $font.Color.RGB = 255
}
